$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 67

$ws.Cells.Item($row, 1).Value = 45748
$ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($row - 1, 1).NumberFormat

$ws.Cells.Item($row, 2).Value = 0
$ws.Cells.Item($row, 3).Value = 0.5
$ws.Cells.Item($row, 4).Value = 1.5
$ws.Cells.Item($row, 5).Value = 3.5
$ws.Cells.Item($row, 6).Value = 4.5
$ws.Cells.Item($row, 7).Value = 9.5
$ws.Cells.Item($row, 8).Value = 14.5
